$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the data value that drives the dependent formulas (H7, I7, H9, I9)
$ws.Range("B7").Value = 164

# Move/update the active selection on the sheet
$ws.Activate()
$ws.Range("F13").Select()
